$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 12/13 — this pushes the existing rows 12..53 down to 14..55,
# carrying their values/styles with them (matches the diff's net effect of a 2-row shift).
$ws.Rows("12:13").Insert()

# Populate the two newly-inserted rows with the new weekly entries.
# Columns A,B,C,E,F,G,H,I,J,T are constant for this market/product across the sheet.

$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C12").Value = "Ñuble"
$ws.Range("D12").Value = [DateTime]"2021-11-25"
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100103
$ws.Range("H12").Value = "Frutos de hueso (carozo)"
$ws.Range("I12").Value = 100103001
$ws.Range("J12").Value = "Cereza"
$ws.Range("K12").Value = "Santina"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 19000
$ws.Range("P12").Value = 18500
$ws.Range("Q12").Value = "$/caja 10 kilos"
$ws.Range("R12").Value = "Provincia de Curicó"
$ws.Range("S12").Value = 1850
$ws.Range("T12").Value = 10

$ws.Range("A13").Value = 7
$ws.Range("B13").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C13").Value = "Ñuble"
$ws.Range("D13").Value = [DateTime]"2021-11-25"
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100103
$ws.Range("H13").Value = "Frutos de hueso (carozo)"
$ws.Range("I13").Value = 100103001
$ws.Range("J13").Value = "Cereza"
$ws.Range("K13").Value = "Santina"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 16000
$ws.Range("O13").Value = 17000
$ws.Range("P13").Value = 16500
$ws.Range("Q13").Value = "$/caja 10 kilos"
$ws.Range("R13").Value = "Provincia de Curicó"
$ws.Range("S13").Value = 1650
$ws.Range("T13").Value = 10
